$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
    "bugredasmissões",
    "C R Juvenal",
    "Doug Leal F.C",
    "Esquadrão Gazembrino",
    "FBC Colorado",
    "GaúchoDaFronteira F.C",
    "GE Bebum",
    "GrioTeam",
    "Grêmio_Campeão_LA_27",
    "JV5 Tricolor Gaúcho",
    "La Primeira Patada Es Nuestra",
    "lsauer fc",
    "Medonho´s F.C.",
    "NHU PORÃ SAF.",
    "Pontaç0 F.C.",
    "SC 100 Sono",
    "SC ÉoINTER!",
    "Texas Club 2026"
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $teams[$i]
}
